# Atualização e criação de arquivos.
# Adds three new certificate rows (91-93) to the "Certificados, cursos, badges"
# sheet, mirroring the formatting of the last existing row (90), and wires up
# the new hyperlink cells in column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Certificados, cursos, badges")

$xlPasteFormats = -4122

# Clone the formatting of the last data row (90) into the three new rows so
# borders / number formats / hyperlink style match the rest of the table.
$ws.Range("B90:I90").Copy($ws.Range("B91:I91"))
$ws.Range("B90:I90").Copy($ws.Range("B92:I92"))
$ws.Range("B90:I90").Copy($ws.Range("B93:I93"))

# Row 91: Data Governance - The Complete Course for Beginners
$ws.Range("B91").Value = "Udemy"
$ws.Range("C91").Value = "Data Governance - The Complete Course for Beginners"
$ws.Range("D91").Value = 5
$ws.Range("E91").Value = 45711
$ws.Range("F91").Value = "https://www.udemy.com/certificate/UC-adc2168e-81f1-45c9-be8a-df6a70772ae0/"
$ws.Range("G91").Value = "Ok"
$ws.Range("H91").Value = "Ok"
$ws.Range("I91").Value = 45711

# Row 92: CDMP Specialist Exam Notes & Practice Questions Data Quality
$ws.Range("B92").Value = "Udemy"
$ws.Range("C92").Value = "CDMP Specialist Exam Notes & Practice Questions Data Quality"
$ws.Range("D92").Value = 3
$ws.Range("E92").Value = 45711
$ws.Range("F92").Value = "https://www.udemy.com/certificate/UC-3d584eb9-7d50-47ac-bcd8-f57bca60910d/"
$ws.Range("G92").Value = "Ok"
$ws.Range("H92").Value = "Ok"
$ws.Range("I92").Value = 45711

# Row 93: Como estabelecer e analisar indicadores de desempenho (KPI)
# (Url entered before the course title, matching the author's original
# shared-string insertion order.)
$ws.Range("B93").Value = "Udemy"
$ws.Range("F93").Value = "https://www.udemy.com/certificate/UC-21f3f67c-238e-47bc-be8a-1a7627df2a87/"
$ws.Range("C93").Value = "Como estabelecer e analisar indicadores de desempenho (KPI)"
$ws.Range("D93").Value = 5
$ws.Range("E93").Value = 45711
$ws.Range("G93").Value = "Ok"
$ws.Range("H93").Value = "Ok"
$ws.Range("I93").Value = 45711

# Wire up the hyperlinks on the new "Url" cells (column F), then restore the
# exact borderered hyperlink cell format (Hyperlinks.Add re-derives the
# builtin Hyperlink style which loses the 9pt font / thin border combo used
# throughout this table).
$ws.Hyperlinks.Add($ws.Range("F91"), "https://www.udemy.com/certificate/UC-adc2168e-81f1-45c9-be8a-df6a70772ae0/")
$ws.Range("F90").Copy()
$ws.Range("F91").PasteSpecial($xlPasteFormats)

$ws.Hyperlinks.Add($ws.Range("F92"), "https://www.udemy.com/certificate/UC-3d584eb9-7d50-47ac-bcd8-f57bca60910d/")
$ws.Range("F90").Copy()
$ws.Range("F92").PasteSpecial($xlPasteFormats)

$ws.Hyperlinks.Add($ws.Range("F93"), "https://www.udemy.com/certificate/UC-21f3f67c-238e-47bc-be8a-1a7627df2a87/")
$ws.Range("F90").Copy()
$ws.Range("F93").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# Match the author's final cursor position / scroll state.
$ws.Range("F93").Select()
$excel.ActiveWindow.ScrollRow = 66
$excel.ActiveWindow.ScrollColumn = 3
